# Regenerate sval data to filter save games.
# Update B2:G8 with the new computed values (TB, d2S, K, IP, Win, sum).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 6.189590430959694)
    3 = @(0.6606524410359556, 0.04071648406533734, 3.537761648806719, 0.4942365360607697, 0, 4.733367109968782)
    4 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 6.189590430959694)
    5 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548)
    6 = @(3.286832544864788, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 0, 4.840633575959121)
    7 = @(0.04271373187048222, 0.306821227259698, 0.1494219747398047, 10.19245300693656, 0, 10.69140994080654)
    8 = @(0.2917716402565462, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 0, 1.242251378316819)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
